$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.442.51"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "'2.102.20"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "'335.49"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5235"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'0.4605"
$ws.Range("E8").Value = "  +7.41%  "
$ws.Range("D9").Value = "'53.26"
$ws.Range("E9").Value = "  +14.56%  "
$ws.Range("D10").Value = "'0.08945"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("D11").Value = "'1.177"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").Value = "'24.42"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "'2.093.55"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "'6.793"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "'7.946"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "'96.38"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "'0.00001133"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "'0.06629"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'19.27"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'30.512.02"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "'2.340.39"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").Value = "'22.30"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "'2.563"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "'163.89"
$ws.Range("D30").Value = "'132.65"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "'1.686"
$ws.Range("E33").Value = "  +10.24%  "
$ws.Range("D34").Value = "'6.155"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").Value = "'3.925"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").Value = "'10.44"
$ws.Range("E36").Value = "  +8.32%  "
$ws.Range("D37").Value = "'0.02571"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'0.06817"
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").Value = "'5.545"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "'12.80"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").Value = "'0.2287"
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").Value = "'0.6879"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("D43").Value = "'1.246"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'2.342"
$ws.Range("E44").Value = "  +6.76%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'14.02"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.6374"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").Value = "'3.664"
$ws.Range("D49").Value = "'0.00000000355"
$ws.Range("E49").Value = "  +24.42%  "
$ws.Range("D50").Value = "'1.245"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "'83.40"
$ws.Range("E51").Value = "  +1.52%  "
